$d = $word.ActiveDocument

function Get-ParaText($para) {
    return $para.Range.Text.TrimEnd([char]13)
}

function Find-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ((Get-ParaText $p) -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Simple in-place text substitutions (do not change paragraph count).
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "New in this update", $true, $false, $false, $false, $false,
    $true, 1, $false, "Document structure update", 2)

$d.Content.Find.Execute(
    "- Added backend store update endpoint: ``PUT /api/stores/{id}`` with owner/admin policy + tenancy checks.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "- Kept ``masterrequirement.docx`` for requirement/status matrix.", 2)

$d.Content.Find.Execute(
    "- Wired Settings General tab save button to backend store update API.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "- Renamed duplicate ``mainrequirementstatus.docx`` to ``release_checklist.docx``.", 2)

$d.Content.Find.Execute(
    "- Added save status feedback on Settings page.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "- ``release_checklist.docx`` now contains strict Owner/ETA/Risk/Priority release checklist.", 2)

$d.Content.Find.Execute(
    "- Last pushed commit: 131e423", $true, $false, $false, $false, $false,
    $true, 1, $false, "- Last pushed commit: 1adb53a", 2)

$d.Content.Find.Execute(
    "- Current store settings API wiring is local and not pushed yet.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "- Current document rename/update is local and not pushed yet.", 2)

# ---------------------------------------------------------------------------
# 2) Structural edits: work from the bottom of the document upward so that
#    paragraph indices referenced below stay valid for the whole script.
# ---------------------------------------------------------------------------

# 2a) Remove the whole "Config" / "DB/Migrations" block: from the "Config"
#     heading paragraph through (but excluding) the "Git state" heading.
$configPara = Find-ParagraphByText $d "Config"
$gitStatePara = Find-ParagraphByText $d "Git state"
if (($configPara -ne $null) -and ($gitStatePara -ne $null)) {
    $blockStart = $configPara.Range.Start
    $blockEnd = $gitStatePara.Range.Start
    $d.Range($blockStart, $blockEnd).Delete()
}

# 2b) Remove the now-trailing 4th bullet under "Document structure update"
#     ("- Backend build verified successful ...").
$buildPara = Find-ParagraphByText $d "- Backend build verified successful (``dotnet build``) with warnings only."
if ($buildPara -ne $null) {
    $buildPara.Range.Delete()
}

# 2c) Insert the new bullet right after "- Product/Customer/Order
#     create+update+delete API wiring done end-to-end." (before the blank
#     paragraph that precedes "Document structure update").
$anchorPara = Find-ParagraphByText $d "- Product/Customer/Order create+update+delete API wiring done end-to-end."
if ($anchorPara -ne $null) {
    $anchorIndex = $anchorPara.Index
    $anchorPara.Range.InsertParagraphAfter()
    $d.Paragraphs.Item($anchorIndex + 1).Range.Text = "- Settings General save now updates store via backend API."
}

Write-Output "done"
